# Add a new worksheet "BeagleBone Received Signals" that maps Beaglebone
# message fields (base mode / system status / light pattern) to the
# Beaglebone Arduino state-machine states, and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "BeagleBone Received Signals"

# Column widths similar to the authored layout (values chosen so the
# engine's internal character-width quantization lands on the closest
# achievable width to the original 17.44140625 / 15.77734375 / 38.21875 /
# 25.6640625 / 31.21875 "characters").
$ws2.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(3).ColumnWidth = 37.333333333333336
$ws2.Columns.Item(4).ColumnWidth = 24.833333333333332
$ws2.Columns.Item(5).ColumnWidth = 30.333333333333332

# Header row (row 2).
$ws2.Range("C2").Value = "base mode (MAV_MODE)"
$ws2.Range("D2").Value = "system status (MAV_STATE)"
$ws2.Range("E2").Value = "Light pattern"

# State definitions in columns A/B.
$ws2.Range("A3").Value = "BONE_POWERUP,"
$ws2.Range("B3").Value = "/**< Initial starting state "

$ws2.Range("A4").Value = "  BONE_SELFTEST,"
$ws2.Range("B4").Value = "/**< Initial self-test"
$ws2.Range("C4").Value = "MAV_MODE_PREFLIGHT"
$ws2.Range("D4").Value = "MAV_STATE_BOOT"

$ws2.Range("A5").Value = "  BONE_DISARMED,"
$ws2.Range("B5").Value = "/**< Disarmed wait state"
$ws2.Range("C5").Value = "MAV_MODE_PREFLIGHT"
$ws2.Range("D5").Value = "MAV_STATE_STANDBY"

$ws2.Range("C6").Value = "MAV_MODE_MANUAL_DISARMED"
$ws2.Range("D6").Value = "MAV_STATE_STANDBY"

$ws2.Range("A7").Value = "  BONE_ARMED,"
$ws2.Range("B7").Value = "/**< Beaglebone armed & ready to navigate"
$ws2.Range("C7").Value = "MAV_MODE_MANUAL_ARMED"
$ws2.Range("D7").Value = "MAV_STATE_STANDBY"

$ws2.Range("C8").Value = "MAV_MODE_AUTO_ARMED"
$ws2.Range("D8").Value = "MAV_STATE_STANDBY"

$ws2.Range("A9").Value = "  BONE_WAYPOINT,"
$ws2.Range("B9").Value = "/**< Beaglebone navigating by waypoints"
$ws2.Range("C9").Value = "MAV_MODE_AUTO_ARMED"
$ws2.Range("D9").Value = "MAV_STATE_ACTIVE"

$ws2.Range("A10").Value = "  BONE_STEERING,"
$ws2.Range("B10").Value = "/**< Beaglebone manual steering"
$ws2.Range("C10").Value = "MAV_MODE_MANUAL_ARMED"
$ws2.Range("D10").Value = "MAV_STATE_ACTIVE"

$ws2.Range("A11").Value = "  BONE_NOSIGNAL,"
$ws2.Range("B11").Value = "/**< Beaglbone has lost shore signal"
$ws2.Range("C11").Value = "MAV_MODE_MANUAL_ARMED"
$ws2.Range("D11").Value = "MAV_STATE_CRITICAL"

$ws2.Range("C12").Value = "MAV_MODE_AUTO_ARMED"
$ws2.Range("D12").Value = "MAV_STATE_CRITICAL"

$ws2.Range("A13").Value = "  BONE_FAULT"
$ws2.Range("B13").Value = "/**< Beaglebone faulted "
$ws2.Range("C13").Value = "Any"
$ws2.Range("D13").Value = "MAV_STATE_EMERGENCY"

# Original sheet keeps its existing selection.
$ws1.Range("A32").Select()

# New sheet becomes the active / selected tab, zoomed to 90%.
$ws2.Activate()
$ws2.Range("B25").Select()
$excel.ActiveWindow.Zoom = 90
